$d = $word.ActiveDocument

# Locate the paragraph that contains the sentence being edited.
$find = $d.Content.Find
$find.ClearFormatting()
$ok = $find.Execute("configuré sur votre bureau.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $ok) {
    throw "Could not locate the target sentence in the document."
}

$targetParagraph = $find.Parent.Paragraphs(1)
$paraRange = $targetParagraph.Range

# Rebuild the whole paragraph, keeping every run untouched except the final
# one, which is split so the (mis-spelled) word "burea" sits in its own run
# flanked by spell-check markers, and the trailing sentence is replaced.
$newParagraphXml = '<w:p w14:paraId="47D156EC" w14:textId="6C458120" w:rsidR="000A4AD1" w:rsidRDefault="001C628F" w:rsidP="00A94B6C">' +
    '<w:r><w:t>Lors de la consultation d’un document dans l’ancienne archive (</w:t></w:r>' +
    '<w:r w:rsidRPr="004643E6"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ELO</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>' +
    '<w:r w:rsidR="004643E6"><w:t>,</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> il est nécessaire de sortir le document en le glissant dans le dossier </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r w:rsidR="004643E6"><w:t>ELO_Transfer</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidR="00B74ED6"><w:t xml:space="preserve"> configuré sur votre </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>burea</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> et de le supprimer de l’archive</w:t></w:r>' +
    '</w:p>'

$paraRange.InsertXML($newParagraphXml)
